$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 143, shifting existing rows 143-188 down to 144-189
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new weekly record
$ws.Cells.Item(143, 1).Value = 9
$ws.Cells.Item(143, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(143, 3).Value = "Metropolitana"
$ws.Cells.Item(143, 4).Value = 44524
$ws.Cells.Item(143, 5).Value = 13
$ws.Cells.Item(143, 6).Value = 300000001
$ws.Cells.Item(143, 7).Value = "Rabanito"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 7900
$ws.Cells.Item(143, 11).Value = 2500
$ws.Cells.Item(143, 12).Value = 3000
$ws.Cells.Item(143, 13).Value = 2747
$ws.Cells.Item(143, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(143, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(143, 16).Value = 27
$ws.Cells.Item(143, 17).Value = 100
$ws.Cells.Item(143, 18).Value = "Hortaliza"
